# Word COM-interop script applying the tracked changes:
#  1. Merge split runs / drop proofErr-bracketed runs for three paragraphs
#     (the visible text is unchanged, only the run/proofErr structure is
#     collapsed into a single run per paragraph).
#  2. Add even/default/first headers & footers to the one section, and set
#     their text (footer2 gets the new "M varshith 192211780" signature,
#     all headers and the even/first footers stay blank).

$d = $word.ActiveDocument

# --- 1. Collapse the three split-run paragraphs into single runs ---------

$d.Content.Find.Execute(
    "To Create a Cloning of a VM and Test it by loading the Previous Version/Cloned  VM.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "To Create a Cloning of a VM and Test it by loading the Previous Version/Cloned  VM.",
    2) | Out-Null

$d.Content.Find.Execute(
    "Then a message will be displayed. In that select the first option which will be : the current state in the VM" + [char]0x201D,
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Then a message will be displayed. In that select the first option which will be : the current state in the VM" + [char]0x201D,
    2) | Out-Null

$d.Content.Find.Execute(
    "then we need to provide the name of the new cloned VM.and then click " + [char]0x201C + "finish" + [char]0x201D,
    $true, $false, $false, $false, $false, $true, 1, $false,
    "then we need to provide the name of the new cloned VM.and then click " + [char]0x201C + "finish" + [char]0x201D,
    2) | Out-Null

# --- 2. Headers / footers --------------------------------------------------

$sec = $d.Sections(1)

# Word's Headers/Footers(1..3) index order is Primary, FirstPage, EvenPages
# (wdHeaderFooterPrimary / wdHeaderFooterFirstPage / wdHeaderFooterEvenPages)
# which this engine persists as header2.xml/header3.xml/header1.xml and
# footer2.xml/footer3.xml/footer1.xml respectively (and wires the matching
# default/first/even w:headerReference+w:footerReference into sectPr).
# Touching all three of each materialises every part; only the primary
# footer (-> footer2.xml) gets real text, everything else stays blank.
$sec.Headers(1).Range.Text = ""
$sec.Headers(2).Range.Text = ""
$sec.Headers(3).Range.Text = ""

$sec.Footers(1).Range.Text = "M varshith 192211780"
$sec.Footers(2).Range.Text = ""
$sec.Footers(3).Range.Text = ""
